# Apply the latest coinranking snapshot values captured by the scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: Address -> new display value.
# Values that look like plain numbers are prefixed with a leading single
# quote so Excel stores them as text (matching the sheets existing
# formatting of Price/Volume columns as inline strings, e.g. "1.00", "0.540").
$updates = [ordered]@{
    "D2" = "42.609.24"
    "E2" = "  -0.04%  "
    "D3" = "2.510.55"
    "E3" = "  -1.53%  "
    "D4" = "'1.00"
    "E4" = "  +0.13%  "
    "D5" = "'315.69"
    "E5" = "  +4.62%  "
    "D6" = "'96.41"
    "E6" = "  -1.05%  "
    "E7" = "  +2.25%  "
    "E8" = "  -0.05%  "
    "D9" = "'0.540"
    "E9" = "  -0.75%  "
    "D10" = "'36.31"
    "E10" = "  +1.03%  "
    "D11" = "'0.0815"
    "E11" = "  +1.22%  "
    "D12" = "'7.70"
    "E12" = "  +2.90%  "
    "E13" = "  -2.68%  "
    "D14" = "2.895.15"
    "E14" = "  -1.53%  "
    "D15" = "'15.59"
    "E15" = "  +6.57%  "
    "D16" = "2.497.59"
    "E16" = "  -1.89%  "
    "D17" = "'0.861"
    "D18" = "42.571.19"
    "E18" = "  -0.23%  "
    "D19" = "'12.92"
    "E19" = "  -2.92%  "
    "D20" = "0.0₃0974"
    "E20" = "  -0.66%  "
    "D21" = "'6.55"
    "E21" = "  +0.00%  "
    "D22" = "'71.61"
    "E22" = "  +0.17%  "
    "D23" = "'253.09"
    "E23" = "  -0.31%  "
    "D24" = "'2.98"
    "E24" = "  +1.59%  "
    "D25" = "'2.05"
    "E25" = "  -1.07%  "
    "D26" = "'27.10"
    "E26" = "  -1.87%  "
    "D27" = "'0.999"
    "E27" = "  -0.13%  "
    "D28" = "'2.35"
    "E28" = "  +12.63%  "
    "D29" = "'10.17"
    "E29" = "  +1.73%  "
    "D30" = "'38.05"
    "E30" = "  +0.86%  "
    "D31" = "'5.93"
    "E31" = "  -0.78%  "
    "D32" = "'157.21"
    "E32" = "  +1.48%  "
    "D33" = "'19.23"
    "E33" = "  +5.68%  "
    "D34" = "'3.32"
    "E34" = "  +1.10%  "
    "D35" = "'0.0789"
    "E35" = "  -1.30%  "
    "D36" = "'2.08"
    "E36" = "  -4.29%  "
    "E37" = "  -4.21%  "
    "E38" = "  -0.54%  "
    "E39" = "  +1.52%  "
    "D40" = "'24.20"
    "E40" = "  -5.47%  "
    "D41" = "'3.41"
    "E41" = "  +1.73%  "
    "D42" = "'3.87"
    "E42" = "  +0.33%  "
    "B43" = "VeChain"
    "C43" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "D43" = "'0.0304"
    "E43" = "  +0.56%  "
    "B44" = "ApeXProtocol"
    "C44" = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
    "D44" = "'2.02"
    "E44" = "  -2.92%  "
    "E45" = "  +0.02%  "
    "D46" = "2.020.85"
    "E46" = "  -1.70%  "
    "D47" = "'84.49"
    "E47" = "  -3.87%  "
    "D48" = "'9.00"
    "E48" = "  -2.13%  "
    "B49" = "ordi"
    "C49" = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
    "D49" = "'74.76"
    "E49" = "  +0.27%  "
    "B50" = "RocketPoolETH"
    "C50" = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
    "D50" = "2.754.47"
    "E50" = "  -1.52%  "
    "E51" = "  +0.65%  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
